$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row: "_old" -> "_FV2404", "_new" -> "_FV2410"
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value()
    if ($text -like "*_old") {
        $cell.Value = ($text -replace "_old$", "_FV2404")
    } elseif ($text -like "*_new") {
        $cell.Value = ($text -replace "_new$", "_FV2410")
    }
}

# 2) Turn the data range into an Excel Table (ListObject) with a header row,
#    without letting Excel bake the existing header formatting into a new
#    dxf (stash it in a scratch row, build the table against a blank
#    header, then restore the formatting from the stash).
$headerRange = $ws.Range("A1:U1")
$tempRange = $ws.Range("A200:U200")

$headerRange.Copy()
$tempRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U79"), [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

$tempRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$tempRange.Clear()

# 3) Freeze the header row (split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
